$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values so they are stored as text, matching source data
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '47.422.17'
$ws.Range("E2").Value = '  +4.66%  '
$ws.Range("D3").Value = '2.498.23'
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '324.07'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '107.96'
$ws.Range("E6").Value = '  +5.25%  '
$ws.Range("E7").Value = '  +2.20%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +3.30%  '
$ws.Range("D10").Value = '38.18'
$ws.Range("E10").Value = '  +7.36%  '
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("D13").Value = '18.46'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").Value = '7.22'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").Value = '2.888.97'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '2.499.69'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = '0.853'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '47.342.86'
$ws.Range("E18").Value = '  +4.68%  '
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").Value = '  +6.24%  '
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = '70.76'
$ws.Range("E22").Value = '  +2.57%  '
$ws.Range("E23").Value = '  +7.29%  '
$ws.Range("D24").Value = '251.41'
$ws.Range("E24").Value = '  +2.68%  '
$ws.Range("E25").Value = '  +4.22%  '
$ws.Range("D26").Value = '26.27'
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '10.09'
$ws.Range("E28").Value = '  +4.80%  '
$ws.Range("D29").Value = '35.39'
$ws.Range("E29").Value = '  +7.14%  '
$ws.Range("D30").Value = '0.138'
$ws.Range("E30").Value = '  +9.35%  '
$ws.Range("E31").Value = '  -8.59%  '
$ws.Range("D32").Value = '49.47'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.53'
$ws.Range("E33").Value = '  +6.06%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = '19.82'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("D37").Value = '2.00'
$ws.Range("E37").Value = '  +6.67%  '
$ws.Range("D38").Value = '4.72'
$ws.Range("E38").Value = '  +6.20%  '
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  +4.41%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").Value = '121.69'
$ws.Range("E42").Value = '  -3.89%  '
$ws.Range("D43").Value = '21.05'
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").Value = '0.0299'
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("D45").Value = '1.970.07'
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").Value = '9.07'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("E50").Value = '  +10.07%  '
$ws.Range("D51").Value = '79.90'
$ws.Range("E51").Value = '  +4.24%  '
